$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 ---
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AG2").Value = 452
$ws.Range("AJ2").Value = 16000000

# --- Row 3 ---
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("D3").Value = 3096
$ws.Range("E3").Value = 1118
$ws.Range("F3").Value = 1118
$ws.Range("G3").Value = 1151
$ws.Range("H3").Value = 837
$ws.Range("I3").Value = 837
$ws.Range("K3").Value = 3880
$ws.Range("L3").Value = 712
$ws.Range("M3").Value = 3168
$ws.Range("N3").Value = 3168
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = 693
$ws.Range("R3").Value = -1878
$ws.Range("S3").Value = 1560
$ws.Range("T3").Value = 56
$ws.Range("U3").Value = 637
$ws.Range("W3").Value = 36.1
$ws.Range("X3").Value = 27.03
$ws.Range("AA3").Value = 22.48
$ws.Range("AB3").Value = 7153.09
$ws.Range("AC3").Value = 5206
$ws.Range("AD3").Value = 14.7
$ws.Range("AE3").Value = 18133
$ws.Range("AF3").Value = 4.22
$ws.Range("AG3").Value = 958
$ws.Range("AH3").Value = 1.25
$ws.Range("AI3").Value = 19.99
$ws.Range("AJ3").Value = 17471300

# --- Row 4 ---
$ws.Range("D4").Value = 3261
$ws.Range("E4").Value = 909
$ws.Range("F4").Value = 733
$ws.Range("G4").Value = 990
$ws.Range("H4").Value = 667
$ws.Range("I4").Value = 355
$ws.Range("J4").Value = 312
$ws.Range("K4").Value = 5394
$ws.Range("L4").Value = 1312
$ws.Range("M4").Value = 4081
$ws.Range("N4").Value = 2157
$ws.Range("O4").Value = 1924
$ws.Range("P4").Value = 46
$ws.Range("Q4").Value = 414
$ws.Range("R4").Value = -138
$ws.Range("S4").Value = -607
$ws.Range("T4").Value = 150
$ws.Range("U4").Value = 264
$ws.Range("V4").Value = 351
$ws.Range("W4").Value = 27.87
$ws.Range("X4").Value = 20.45
$ws.Range("Y4").Value = 13.34
$ws.Range("Z4").Value = 14.38
$ws.Range("AA4").Value = 32.16
$ws.Range("AB4").Value = 4515.79
$ws.Range("AC4").Value = 2033
$ws.Range("AD4").Value = 21.15
$ws.Range("AE4").Value = 12601
$ws.Range("AF4").Value = 3.41
$ws.Range("AG4").Value = 701
$ws.Range("AH4").Value = 1.63
$ws.Range("AI4").Value = 33.8
$ws.Range("AJ4").Value = 17471300

# --- Row 5 ---
$ws.Range("D5").Value = 2457
$ws.Range("E5").Value = 454
$ws.Range("F5").Value = 454
$ws.Range("G5").Value = 488
$ws.Range("H5").Value = 420
$ws.Range("I5").Value = 275
$ws.Range("J5").Value = 145
$ws.Range("K5").Value = 5624
$ws.Range("L5").Value = 931
$ws.Range("M5").Value = 4693
$ws.Range("N5").Value = 4373
$ws.Range("O5").Value = 320
$ws.Range("P5").Value = 110
$ws.Range("Q5").Value = 285
$ws.Range("R5").Value = -352
$ws.Range("S5").Value = -88
$ws.Range("T5").Value = 122
$ws.Range("U5").Value = 163
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 18.46
$ws.Range("X5").Value = 17.1
$ws.Range("Y5").Value = 8.42
$ws.Range("Z5").Value = 7.63
$ws.Range("AA5").Value = 19.84
$ws.Range("AB5").Value = 3897.37
$ws.Range("AC5").Value = 1345
$ws.Range("AD5").Value = 33.76
$ws.Range("AE5").Value = 24821
$ws.Range("AF5").Value = 1.83
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 1.32
$ws.Range("AI5").Value = 38.44
$ws.Range("AJ5").Value = 21929315

# --- Row 6 ---
$ws.Range("D6").Value = 2141
$ws.Range("E6").Value = 199
$ws.Range("F6").Value = 199
$ws.Range("G6").Value = 260
$ws.Range("H6").Value = 186
$ws.Range("I6").Value = 82
$ws.Range("K6").Value = 5356
$ws.Range("L6").Value = 608
$ws.Range("M6").Value = 4747
$ws.Range("N6").Value = 4349
$ws.Range("P6").Value = 110
$ws.Range("Q6").Value = 39
$ws.Range("R6").Value = -11
$ws.Range("S6").Value = -301
$ws.Range("T6").Value = 69
$ws.Range("U6").Value = -30
$ws.Range("V6").Value = 65
$ws.Range("W6").Value = 9.300000000000001
$ws.Range("X6").Value = 8.67
$ws.Range("Y6").Value = 1.89
$ws.Range("Z6").Value = 3.38
$ws.Range("AA6").Value = 12.81
$ws.Range("AB6").Value = 3873.89
$ws.Range("AC6").Value = 376
$ws.Range("AD6").Value = 71.45999999999999
$ws.Range("AE6").Value = 24688
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 0.5600000000000001
$ws.Range("AI6").Value = 32.07
$ws.Range("AJ6").Value = 21929315

# --- Row 7 ---
$ws.Range("U7").ClearContents()
$ws.Range("D7").Value = 2096
$ws.Range("E7").Value = 107
$ws.Range("G7").Value = 173
$ws.Range("H7").Value = 47
$ws.Range("I7").Value = -29
$ws.Range("K7").Value = 5578
$ws.Range("L7").Value = 813
$ws.Range("M7").Value = 4765
$ws.Range("N7").Value = 4291
$ws.Range("P7").Value = 110
$ws.Range("Q7").Value = 168
$ws.Range("R7").Value = -282
$ws.Range("S7").Value = 209
$ws.Range("T7").Value = 90
$ws.Range("W7").Value = 5.1
$ws.Range("X7").Value = 2.24
$ws.Range("Y7").Value = -0.67
$ws.Range("Z7").Value = 0.86
$ws.Range("AA7").Value = 17.06
$ws.Range("AC7").Value = -132
$ws.Range("AD7").Value = -133.47
$ws.Range("AE7").Value = 24358
$ws.Range("AF7").Value = 0.72
$ws.Range("AG7").Value = 150
$ws.Range("AH7").Value = 0.85
$ws.Range("AI7").Value = -113.43

# --- Row 8 ---
$ws.Range("U8").ClearContents()
$ws.Range("D8").Value = 2269
$ws.Range("E8").Value = 255
$ws.Range("G8").Value = 328
$ws.Range("H8").Value = 249
$ws.Range("I8").Value = 149
$ws.Range("K8").Value = 5650
$ws.Range("L8").Value = 666
$ws.Range("M8").Value = 4984
$ws.Range("N8").Value = 4410
$ws.Range("P8").Value = 110
$ws.Range("Q8").Value = 161
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = -126
$ws.Range("T8").Value = 120
$ws.Range("W8").Value = 11.24
$ws.Range("X8").Value = 10.97
$ws.Range("Y8").Value = 3.43
$ws.Range("Z8").Value = 4.44
$ws.Range("AA8").Value = 13.36
$ws.Range("AC8").Value = 679
$ws.Range("AD8").Value = 25.98
$ws.Range("AE8").Value = 25033
$ws.Range("AF8").Value = 0.71
$ws.Range("AG8").Value = 150
$ws.Range("AH8").Value = 0.85
$ws.Range("AI8").Value = 22.08

# --- Row 9 ---
$ws.Range("U9").ClearContents()
$ws.Range("D9").Value = 2503
$ws.Range("E9").Value = 310
$ws.Range("G9").Value = 380
$ws.Range("H9").Value = 289
$ws.Range("I9").Value = 173
$ws.Range("K9").Value = 5874
$ws.Range("L9").Value = 631
$ws.Range("M9").Value = 5243
$ws.Range("N9").Value = 4554
$ws.Range("P9").Value = 110
$ws.Range("Q9").Value = 291
$ws.Range("R9").Value = -5
$ws.Range("S9").Value = -26
$ws.Range("T9").Value = 120
$ws.Range("W9").Value = 12.38
$ws.Range("X9").Value = 11.55
$ws.Range("Y9").Value = 3.86
$ws.Range("Z9").Value = 5.02
$ws.Range("AA9").Value = 12.04
$ws.Range("AC9").Value = 789
$ws.Range("AD9").Value = 22.37
$ws.Range("AF9").Value = 0.68
$ws.Range("AG9").Value = 150
$ws.Range("AH9").Value = 0.85
$ws.Range("AI9").Value = 19.01
